# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = @'
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 3.7 = 14177.45 pesos
✅ 14177.45 pesos = 3.68 = 963.55 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
'@
$ws1.Range("A1").Value = $newText

# --- tasas: update the rate cells used by the conversion formulas ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 270.5
$ws2.Range("O10").Value = 3835
$ws2.Range("N12").Value = 3855
$ws2.Range("O12").Value = 262
